# Fixed Bento 80 Test scripts
#
# The "startup" sheet drives the TC01_INS_Filter_Prog-CCDI test case. Cell
# B2 holds the Cypher query used to pull the CCDI project list from Neo4j;
# it is being updated to sort the results and cap the row count, and the
# sheet's active selection is being moved onto that query cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$query = @'
MATCH (p:project)-->(pgm:program)
where pgm.program_id='CCDI'
WITH DISTINCT p, pgm
RETURN
coalesce(p.project_id, '') AS `Project ID`,
coalesce (pgm.program_id, '')AS `Program`,
coalesce(p.project_title, '') AS `Project Title`,
coalesce(p.principal_investigators, '') AS `Principal Investigators`,
coalesce(p.program_officers, '') AS `Program Officers`,
coalesce(p.lead_doc, '')AS `Lead DOC`,
SUBSTRING(p.project_id, 1, 3) AS `Activity code`,
"$" + apoc.number.format(toInteger(p.award_amount)) AS `Award Amount`,
coalesce(p.project_end_date, '') AS `Project End Date`,
coalesce(p.fiscal_year,'')AS `Fiscal Year` Order By p.project_id Asc LIMIT 100
'@

$ws.Range("B2").Value2 = $query

# Move the sheet's selection from the old C15 anchor onto the edited query
# cell (B2), matching the new <selection activeCell="B2" sqref="B2"/>.
[void]$ws.Activate()
[void]$ws.Range("B2").Select()
